# Daily attendance processing - rotate the "Recorded By" (column G) list
# left by one position for every row whose value contains more than one
# comma-separated entry (e.g. "a, b, c" -> "b, c, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)
    $text = $cell.Value2

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = $rotated -join ", "
        }
    }
}
